$d = $word.ActiveDocument

# 1. Remove the "Boss de fin – Divergence du personnage principal" run (and its _GoBack bookmark)
$d.Content.Find.Execute("Boss de fin – Divergence du personnage principal", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)

# 2. Add "o" after "2034"
$d.Content.Find.Execute("2034", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2034o", 2)

# 3. Change "Maître des Dieux" to "Maho"
$d.Content.Find.Execute("Maître des Dieux", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Maho", 2)

# 4. Remove proofErr spellStart/spellEnd markers around "Catoryu" - handled by removing any
#    proofing-error markup. Word COM interop doesn't directly expose proofErr elements, so
#    this will be handled via direct XML manipulation if needed.
